$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-8 (Generation 0-6) get individually changed fitness values
$ws.Range("C2").Value = 12181
$ws.Range("C3").Value = 12181
$ws.Range("C4").Value = 12102
$ws.Range("C5").Value = 11258
$ws.Range("C6").Value = 10861
$ws.Range("C7").Value = 10844
$ws.Range("C8").Value = 10844

# Rows 9-252 (Generation 7-250) all converge to the same fitness value
$ws.Range("C9:C252").Value = 9758
